$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy each year-2020 (column Q) cell format into the new year-2021 (column R) cell,
# then set the 2021 values - mirrors how the sheet author extended the table with a new column.
$ws.Range("Q3").Copy()
$ws.Range("R3").PasteSpecial(-4122)
$ws.Range("R3").Value = 2021

$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Value = 0.12641839647678207

$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
$ws.Range("R5").Value = 0.14922981985616976

$ws.Range("Q6").Copy()
$ws.Range("R6").PasteSpecial(-4122)
$ws.Range("R6").Value = 0.10326895933792253

$ws.Range("Q7").Copy()
$ws.Range("R7").PasteSpecial(-4122)
$ws.Range("R7").Value = 3.433011112114915E-2

$ws.Range("Q8").Copy()
$ws.Range("R8").PasteSpecial(-4122)
$ws.Range("R8").Value = 3.6820478077087354E-2

$ws.Range("Q9").Copy()
$ws.Range("R9").PasteSpecial(-4122)
$ws.Range("R9").Value = 3.1930519190242035E-2

$ws.Range("Q10").Copy()
$ws.Range("R10").PasteSpecial(-4122)
$ws.Range("R10").Value = 8.7302929367211068E-2

$ws.Range("Q11").Copy()
$ws.Range("R11").PasteSpecial(-4122)
$ws.Range("R11").Value = 0.10296328329317765

$ws.Range("Q12").Copy()
$ws.Range("R12").PasteSpecial(-4122)
$ws.Range("R12").Value = 7.1859056271889668E-2

$ws.Range("Q13").Copy()
$ws.Range("R13").PasteSpecial(-4122)
$ws.Range("R13").Value = 0.10716050460690947

$ws.Range("Q14").Copy()
$ws.Range("R14").PasteSpecial(-4122)
$ws.Range("R14").Value = 7.9035451351703812E-2

$ws.Range("Q15").Copy()
$ws.Range("R15").PasteSpecial(-4122)
$ws.Range("R15").Value = 0.13553052227085377

$ws.Range("Q16").Copy()
$ws.Range("R16").PasteSpecial(-4122)
$ws.Range("R16").Value = 6.479643687803946E-2

$ws.Range("Q17").Copy()
$ws.Range("R17").PasteSpecial(-4122)
$ws.Range("R17").Value = 7.643825526207898E-2

$ws.Range("Q18").Copy()
$ws.Range("R18").PasteSpecial(-4122)
$ws.Range("R18").Value = 5.3576570965516782E-2

$ws.Range("Q19").Copy()
$ws.Range("R19").PasteSpecial(-4122)
$ws.Range("R19").Value = 5.4163459619715498E-2

$ws.Range("Q20").Copy()
$ws.Range("R20").PasteSpecial(-4122)
$ws.Range("R20").Value = 6.4872252119520635E-2

$ws.Range("Q21").Copy()
$ws.Range("R21").PasteSpecial(-4122)
$ws.Range("R21").Value = 4.3693418784505472E-2

$ws.Range("Q22").Copy()
$ws.Range("R22").PasteSpecial(-4122)
$ws.Range("R22").Value = 5.1373884452794741E-2

$ws.Range("Q23").Copy()
$ws.Range("R23").PasteSpecial(-4122)
$ws.Range("R23").Value = 2.9662368095156877E-2

$ws.Range("Q24").Copy()
$ws.Range("R24").PasteSpecial(-4122)
$ws.Range("R24").Value = 7.2642215296997686E-2

$ws.Range("Q25").Copy()
$ws.Range("R25").PasteSpecial(-4122)
$ws.Range("R25").Value = 0.13772601093442507

$ws.Range("Q26").Copy()
$ws.Range("R26").PasteSpecial(-4122)
$ws.Range("R26").Value = 0.15668565643254884

$ws.Range("Q27").Copy()
$ws.Range("R27").PasteSpecial(-4122)
$ws.Range("R27").Value = 0.11816042869432726

$ws.Range("Q28").Copy()
$ws.Range("R28").PasteSpecial(-4122)
$ws.Range("R28").Value = 0.33417383115107696

$ws.Range("Q29").Copy()
$ws.Range("R29").PasteSpecial(-4122)
$ws.Range("R29").Value = 0.41139191068108794

$ws.Range("Q30").Copy()
$ws.Range("R30").PasteSpecial(-4122)
$ws.Range("R30").Value = 0.24697746624641295

$ws.Range("Q31").Copy()
$ws.Range("R31").PasteSpecial(-4122)
$ws.Range("R31").Value = 0.16773611144997194

$ws.Range("Q32").Copy()
$ws.Range("R32").PasteSpecial(-4122)
$ws.Range("R32").Value = 0.1959922553363346

$ws.Range("Q33").Copy()
$ws.Range("R33").PasteSpecial(-4122)
$ws.Range("R33").Value = 0.13791201213625709

$ws.Range("Q34").Copy()
$ws.Range("R34").PasteSpecial(-4122)

$ws.Range("Q35").Copy()
$ws.Range("R35").PasteSpecial(-4122)
$ws.Range("R35").Value = 0

$ws.Range("Q36").Copy()
$ws.Range("R36").PasteSpecial(-4122)
$ws.Range("R36").Value = 0.1

$ws.Range("Q37").Copy()
$ws.Range("R37").PasteSpecial(-4122)
$ws.Range("R37").Value = 0.2

$excel.CutCopyMode = $false

# Move the active selection to C1, matching the saved view state.
$ws.Range("C1").Select()